$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '67.061.03'
Set-TextValue 2 5 '  +5.04%  '

# Row 3
Set-TextValue 3 4 '3.467.82'
Set-TextValue 3 5 '  +5.68%  '

# Row 4
Set-TextValue 4 5 '  -0.12%  '

# Row 5
Set-TextValue 5 4 '185.95'
Set-TextValue 5 5 '  +6.05%  '

# Row 6
Set-TextValue 6 4 '545.74'
Set-TextValue 6 5 '  +4.46%  '

# Row 7
Set-TextValue 7 4 '0.612'
Set-TextValue 7 5 '  +1.54%  '

# Row 8
Set-TextValue 8 4 '3.474.09'
Set-TextValue 8 5 '  +6.24%  '

# Row 9
Set-TextValue 9 4 '0.999'
Set-TextValue 9 5 '  -0.38%  '

# Row 10
Set-TextValue 10 4 '0.637'
Set-TextValue 10 5 '  +5.59%  '

# Row 11
Set-TextValue 11 4 '55.92'
Set-TextValue 11 5 '  -1.62%  '

# Row 12
Set-TextValue 12 5 '  +11.50%  '

# Row 13
Set-TextValue 13 4 '0.0000273'
Set-TextValue 13 5 '  +6.32%  '

# Row 14
Set-TextValue 14 4 '9.45'
Set-TextValue 14 5 '  +4.79%  '

# Row 15
Set-TextValue 15 4 '4.015.40'
Set-TextValue 15 5 '  +5.27%  '

# Row 16
Set-TextValue 16 4 '3.465.96'
Set-TextValue 16 5 '  +5.46%  '

# Row 17
Set-TextValue 17 4 '67.388.73'
Set-TextValue 17 5 '  +5.59%  '

# Row 18
Set-TextValue 18 5 '  +3.35%  '

# Row 19
Set-TextValue 19 4 '18.27'
Set-TextValue 19 5 '  +5.47%  '

# Row 20
Set-TextValue 20 4 '11.80'
Set-TextValue 20 5 '  +7.21%  '

# Row 21
Set-TextValue 21 5 '  +6.06%  '

# Row 22
Set-TextValue 22 4 '405.43'
Set-TextValue 22 5 '  +9.05%  '

# Row 23
Set-TextValue 23 5 '  +8.64%  '

# Row 24
Set-TextValue 24 4 '84.69'
Set-TextValue 24 5 '  +5.59%  '

# Row 25
Set-TextValue 25 4 '3.88'
Set-TextValue 25 5 '  +3.62%  '

# Row 26
Set-TextValue 26 4 '4.18'
Set-TextValue 26 5 '  +9.03%  '

# Row 27
Set-TextValue 27 4 '2.91'
Set-TextValue 27 5 '  +10.07%  '

# Row 28
Set-TextValue 28 4 '6.10'
Set-TextValue 28 5 '  -0.85%  '

# Row 29
Set-TextValue 29 4 '11.74'
Set-TextValue 29 5 '  +3.84%  '

# Row 30
Set-TextValue 30 5 '  +4.42%  '

# Row 31
Set-TextValue 31 4 '30.15'
Set-TextValue 31 5 '  +5.24%  '

# Row 32
Set-TextValue 32 4 '675.47'
Set-TextValue 32 5 '  +6.46%  '

# Row 33
Set-TextValue 33 4 '6.85'
Set-TextValue 33 5 '  +4.16%  '

# Row 34
Set-TextValue 34 4 '11.65'
Set-TextValue 34 5 '  +4.25%  '

# Row 35
Set-TextValue 35 5 '  +5.27%  '

# Row 36
Set-TextValue 36 4 '59.10'
Set-TextValue 36 5 '  +0.51%  '

# Row 37
Set-TextValue 37 4 '0.0₃0823'
Set-TextValue 37 5 '  +17.74%  '

# Row 38
Set-TextValue 38 4 '38.55'
Set-TextValue 38 5 '  +6.09%  '

# Row 39
Set-TextValue 39 5 '  +4.39%  '

# Row 40
Set-TextValue 40 5 '  -0.01%  '

# Row 41
Set-TextValue 41 2 'Stacks'
Set-TextValue 41 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 41 4 '3.37'
Set-TextValue 41 5 '  +21.96%  '

# Row 42
Set-TextValue 42 2 'Fetch.AI'
Set-TextValue 42 3 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 42 4 '2.80'
Set-TextValue 42 5 '  +14.76%  '

# Row 43
Set-TextValue 43 5 '  +7.31%  '

# Row 44
Set-TextValue 44 5 '  +0.25%  '

# Row 45
Set-TextValue 45 2 'ThetaToken'
Set-TextValue 45 3 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 45 4 '2.99'
Set-TextValue 45 5 '  +11.79%  '

# Row 46
Set-TextValue 46 2 'Maker'
Set-TextValue 46 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 46 4 '3.045.43'
Set-TextValue 46 5 '  +4.66%  '

# Row 47
Set-TextValue 47 4 '0.0418'
Set-TextValue 47 5 '  +5.97%  '

# Row 48
Set-TextValue 48 4 '3.27'
Set-TextValue 48 5 '  +9.43%  '

# Row 49
Set-TextValue 49 4 '2.74'
Set-TextValue 49 5 '  +3.49%  '

# Row 50
Set-TextValue 50 4 '8.73'
Set-TextValue 50 5 '  +11.83%  '

# Row 51
Set-TextValue 51 2 'Stellar'
Set-TextValue 51 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 51 4 '0.129'
Set-TextValue 51 5 '  +3.58%  '
